# Update the year embedded in the astromap link, e.g.
#   (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).
# becomes
#   (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).
#
# The old text is split across three differently-formatted runs
# ("(" / the hyperlink-styled URL / ")."). We replace the whole
# parenthesised link in one go so the result collapses into a single
# plain run, matching how Word merges runs when the replacement text
# spans multiple differently-formatted runs.

$d = $word.ActiveDocument

$old = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$new = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    [void]$r.MoveEnd(1, -1)
    if ($r.Text -eq $old) {
        $r.Delete()
        $r.InsertAfter($new)
        $found = $true
    }
}

Write-Output "replaced=$found"
